# update scripts wuth new tpm
# Refreshes the NATMI ligand/receptor/edge expression & specificity metrics
# (columns G-T) for the Lama5-Bcam LR-pair sheet with the recomputed TPM-based
# values, leaving the categorical/count columns (A-F, K, L) untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.215324
$ws.Range("H2").Value = 60.645972
$ws.Range("I2").Value = 0.5046434147095277
$ws.Range("J2").Value = 0.5046434147095277
$ws.Range("M2").Value = 44.22587833333333
$ws.Range("N2").Value = 132.677635
$ws.Range("O2").Value = 0.3992109918068397
$ws.Range("P2").Value = 0.3992109918068397
$ws.Range("Q2").Value = 894.0404596929133
$ws.Range("R2").Value = 8046.364137236221
$ws.Range("S2").Value = 0.2014591980949808
$ws.Range("T2").Value = 0.2014591980949809
$ws.Range("G3").Value = 20.215324
$ws.Range("H3").Value = 60.645972
$ws.Range("I3").Value = 0.5046434147095277
$ws.Range("J3").Value = 0.5046434147095277
$ws.Range("O3").Value = 0.007198448896233843
$ws.Range("P3").Value = 0.007198448896233844
$ws.Range("Q3").Value = 16.12106052274933
$ws.Range("R3").Value = 145.089544704744
$ws.Range("S3").Value = 0.003632649831607477
$ws.Range("T3").Value = 0.003632649831607477
$ws.Range("G4").Value = 20.215324
$ws.Range("H4").Value = 60.645972
$ws.Range("I4").Value = 0.5046434147095277
$ws.Range("J4").Value = 0.5046434147095277
$ws.Range("M4").Value = 65.75987233333333
$ws.Range("O4").Value = 0.5935905592969265
$ws.Range("P4").Value = 0.5935905592969265
$ws.Range("Q4").Value = 1329.357125416969
$ws.Range("R4").Value = 11964.21412875272
$ws.Range("S4").Value = 0.2995515667829394
$ws.Range("T4").Value = 0.2995515667829394
$ws.Range("I5").Value = 0.01888864714138046
$ws.Range("J5").Value = 0.01888864714138046
$ws.Range("M5").Value = 44.22587833333333
$ws.Range("N5").Value = 132.677635
$ws.Range("O5").Value = 0.3992109918068397
$ws.Range("P5").Value = 0.3992109918068397
$ws.Range("Q5").Value = 33.46365826051111
$ws.Range("R5").Value = 301.1729243446
$ws.Range("S5").Value = 0.007540555559199921
$ws.Range("T5").Value = 0.00754055555919992
$ws.Range("I6").Value = 0.01888864714138046
$ws.Range("J6").Value = 0.01888864714138046
$ws.Range("O6").Value = 0.007198448896233843
$ws.Range("P6").Value = 0.007198448896233844
$ws.Range("S6").Value = 0.0001359689611662207
$ws.Range("T6").Value = 0.0001359689611662207
$ws.Range("I7").Value = 0.01888864714138046
$ws.Range("J7").Value = 0.01888864714138046
$ws.Range("M7").Value = 65.75987233333333
$ws.Range("O7").Value = 0.5935905592969265
$ws.Range("P7").Value = 0.5935905592969265
$ws.Range("Q7").Value = 49.75742660059111
$ws.Range("S7").Value = 0.01121212262101432
$ws.Range("T7").Value = 0.01121212262101432
$ws.Range("H8").Value = 57.25995900000001
$ws.Range("I8").Value = 0.4764679381490919
$ws.Range("J8").Value = 0.4764679381490919
$ws.Range("M8").Value = 44.22587833333333
$ws.Range("N8").Value = 132.677635
$ws.Range("O8").Value = 0.3992109918068397
$ws.Range("P8").Value = 0.3992109918068397
$ws.Range("Q8").Value = 844.1239933685517
$ws.Range("R8").Value = 7597.115940316967
$ws.Range("S8").Value = 0.1902112381526589
$ws.Range("T8").Value = 0.1902112381526589
$ws.Range("H9").Value = 57.25995900000001
$ws.Range("I9").Value = 0.4764679381490919
$ws.Range("J9").Value = 0.4764679381490919
$ws.Range("O9").Value = 0.007198448896233843
$ws.Range("P9").Value = 0.007198448896233844
$ws.Range("S9").Value = 0.003429830103460145
$ws.Range("T9").Value = 0.003429830103460146
$ws.Range("H10").Value = 57.25995900000001
$ws.Range("I10").Value = 0.4764679381490919
$ws.Range("J10").Value = 0.4764679381490919
$ws.Range("M10").Value = 65.75987233333333
$ws.Range("O10").Value = 0.5935905592969265
$ws.Range("P10").Value = 0.5935905592969265
$ws.Range("S10").Value = 0.2828268698929728
$ws.Range("T10").Value = 0.2828268698929728
